$p = $ppt.ActivePresentation

# Remove the second slide (the bivariate plot slide) from the deck.
$p.Slides.Item(2).Delete()
